# Updating how thermal threshold (column D, "Thresh_temp") is calculated in
# ectopic expression such that the absolute value of the response must be at
# least 3*STD(wt response) above the WT response.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Recalculated Thresh_temp (column D) values ---
$ws.Range("D43").Value = 24
$ws.Range("D44").Value = 21.7
$ws.Range("D46").Value = 40.3
$ws.Range("D47").Value = 31
$ws.Range("D48").Value = 22
$ws.Range("D49").Value = 21.6
$ws.Range("D50").Value = 21.55
$ws.Range("D51").Value = 23.9
$ws.Range("D54").Value = 29.4
$ws.Range("D55").Value = 27.3

# Row 52 (EAH377, 201204_007 / strain XL115) no longer clears the thermal
# threshold: the pair of Thresh_temp / Tmax values for that trial are removed.
$ws.Range("D52:E52").Clear()

# --- View-state bookkeeping to mirror the author's last-saved window state ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 111
$win.ScrollRow = 28
$win.ScrollColumn = 1
$ws.Range("H47").Select()
